# Remove the trailing "Ver no Jupiter..." / copyright paragraphs (and the
# blank paragraph right before them) that follow the requirements line,
# while leaving the blank paragraph + page-break paragraph at the very end
# of the document untouched.

$d = $word.ActiveDocument

# Paragraph.Range.Text includes the trailing paragraph-mark character (CR),
# so strip it before comparing.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "Ver no Jupiter Salvar em pdf Salvar em docx") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $blank = $d.Paragraphs.Item($target - 1)
    $credit = $d.Paragraphs.Item($target + 1)

    $rangeStart = $blank.Range.Start
    $rangeEnd = $credit.Range.End

    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}
